$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standardise the parameter name "cost_variable" -> "cost_variable_om"
# for the block of rows that used the old name (rows 10-39, column C).
$ws.Range("C10:C39").Value = "cost_variable_om"

# Reflect the edited range as the active selection, like Excel would
# leave it after a fill-down / replace operation over that range.
$ws.Range("C10:C39").Select()
